$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value2 = '23.673.99'
$ws.Cells.Item(2, 5).Value2 = '  +1.11%  '
$ws.Cells.Item(3, 4).Value2 = '1.653.22'
$ws.Cells.Item(3, 5).Value2 = '  +1.43%  '
$ws.Cells.Item(4, 4).Value2 = '0.9997'
$ws.Cells.Item(4, 5).Value2 = '  -0.16%  '
$ws.Cells.Item(5, 5).Value2 = '  -0.28%  '
$ws.Cells.Item(6, 4).Value2 = '304.56'
$ws.Cells.Item(6, 5).Value2 = '  +0.21%  '
$ws.Cells.Item(7, 4).Value2 = '0.3818'
$ws.Cells.Item(7, 5).Value2 = '  +1.13%  '
$ws.Cells.Item(8, 4).Value2 = '52.17'
$ws.Cells.Item(8, 5).Value2 = '  +2.35%  '
$ws.Cells.Item(9, 4).Value2 = '0.3619'
$ws.Cells.Item(9, 5).Value2 = '  -0.75%  '
$ws.Cells.Item(10, 4).Value2 = '1.249'
$ws.Cells.Item(10, 5).Value2 = '  +1.30%  '
$ws.Cells.Item(11, 4).Value2 = '0.08222'
$ws.Cells.Item(11, 5).Value2 = '  +0.00%  '
$ws.Cells.Item(12, 4).Value2 = '0.9999'
$ws.Cells.Item(12, 5).Value2 = '  -0.15%  '
$ws.Cells.Item(13, 4).Value2 = '22.57'
$ws.Cells.Item(13, 5).Value2 = '  +0.96%  '
$ws.Cells.Item(14, 4).Value2 = '6.532'
$ws.Cells.Item(14, 5).Value2 = '  -0.09%  '
$ws.Cells.Item(15, 4).Value2 = '7.394'
$ws.Cells.Item(15, 5).Value2 = '  +0.91%  '
$ws.Cells.Item(16, 4).Value2 = '0.00001231'
$ws.Cells.Item(16, 5).Value2 = '  -1.33%  '
$ws.Cells.Item(17, 4).Value2 = '1.656.26'
$ws.Cells.Item(17, 5).Value2 = '  +1.60%  '
$ws.Cells.Item(18, 4).Value2 = '96.74'
$ws.Cells.Item(18, 5).Value2 = '  +2.93%  '
$ws.Cells.Item(19, 4).Value2 = '0.06971'
$ws.Cells.Item(19, 5).Value2 = '  -0.10%  '
$ws.Cells.Item(20, 4).Value2 = '6.791'
$ws.Cells.Item(20, 5).Value2 = '  +4.47%  '
$ws.Cells.Item(21, 4).Value2 = '17.65'
$ws.Cells.Item(21, 5).Value2 = '  -0.12%  '
$ws.Cells.Item(22, 4).Value2 = '0.9987'
$ws.Cells.Item(22, 5).Value2 = '  -0.22%  '
$ws.Cells.Item(23, 4).Value2 = '12.61'
$ws.Cells.Item(23, 5).Value2 = '  -0.63%  '
$ws.Cells.Item(24, 4).Value2 = '23.680.42'
$ws.Cells.Item(24, 5).Value2 = '  +1.11%  '
$ws.Cells.Item(25, 4).Value2 = '2.527'
$ws.Cells.Item(25, 5).Value2 = '  +3.07%  '
$ws.Cells.Item(26, 4).Value2 = '3.094'
$ws.Cells.Item(26, 5).Value2 = '  -0.99%  '
$ws.Cells.Item(27, 4).Value2 = '21.38'
$ws.Cells.Item(27, 5).Value2 = '  +0.08%  '
$ws.Cells.Item(28, 4).Value2 = '152.41'
$ws.Cells.Item(28, 5).Value2 = '  +1.57%  '
$ws.Cells.Item(29, 4).Value2 = '5.203'
$ws.Cells.Item(29, 5).Value2 = '  -1.60%  '
$ws.Cells.Item(30, 4).Value2 = '135.02'
$ws.Cells.Item(30, 5).Value2 = '  +0.96%  '
$ws.Cells.Item(31, 4).Value2 = '1.832.43'
$ws.Cells.Item(31, 5).Value2 = '  +1.26%  '
$ws.Cells.Item(32, 4).Value2 = '6.888'
$ws.Cells.Item(32, 5).Value2 = '  +0.94%  '
$ws.Cells.Item(33, 4).Value2 = '1.088'
$ws.Cells.Item(33, 5).Value2 = '  +5.76%  '
$ws.Cells.Item(34, 4).Value2 = '2.096'
$ws.Cells.Item(34, 5).Value2 = '  -7.36%  '
$ws.Cells.Item(35, 4).Value2 = '11.59'
$ws.Cells.Item(35, 5).Value2 = '  +7.26%  '
$ws.Cells.Item(36, 4).Value2 = '0.02806'
$ws.Cells.Item(36, 5).Value2 = '  +0.86%  '
$ws.Cells.Item(37, 4).Value2 = '0.2517'
$ws.Cells.Item(38, 5).Value2 = '  +2.26%  '
$ws.Cells.Item(39, 5).Value2 = '  +0.78%  '
$ws.Cells.Item(40, 4).Value2 = '0.07048'
$ws.Cells.Item(40, 5).Value2 = '  -0.72%  '
$ws.Cells.Item(41, 4).Value2 = '12.75'
$ws.Cells.Item(41, 5).Value2 = '  +5.24%  '
$ws.Cells.Item(42, 4).Value2 = '0.7065'
$ws.Cells.Item(42, 5).Value2 = '  +0.54%  '
$ws.Cells.Item(43, 4).Value2 = '1.338'
$ws.Cells.Item(43, 5).Value2 = '  -0.56%  '
$ws.Cells.Item(44, 4).Value2 = '16.06'
$ws.Cells.Item(44, 5).Value2 = '  -0.09%  '
$ws.Cells.Item(45, 4).Value2 = '0.6505'
$ws.Cells.Item(45, 5).Value2 = '  -0.54%  '
$ws.Cells.Item(46, 4).Value2 = '2.337'
$ws.Cells.Item(46, 5).Value2 = '  +1.85%  '
$ws.Cells.Item(47, 4).Value2 = '0.9983'
$ws.Cells.Item(47, 5).Value2 = '  -0.19%  '
$ws.Cells.Item(48, 5).Value2 = '  +0.14%  '
$ws.Cells.Item(49, 4).Value2 = '0.08001'
$ws.Cells.Item(49, 5).Value2 = '  -0.19%  '
$ws.Cells.Item(50, 4).Value2 = '128.20'
$ws.Cells.Item(50, 5).Value2 = '  +1.88%  '
$ws.Cells.Item(51, 4).Value2 = '1.192'
$ws.Cells.Item(51, 5).Value2 = '  -0.64%  '
